# cpi.xlsx update: roll the workbook from "through 2018" to "through 2019"
# and append the 2019 CPI-U row pulled from the BLS December 2019 table.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("Data")

# --- Data sheet: append the new 2019 row ------------------------------
# Doing this FIRST (before touching the About-sheet hyperlink text) means
# the shared-string slot that used to hold the old BLS URL gets reused
# for the "2019....." label text, and the new URL string is appended
# fresh afterwards — matching how the shared-string table grows.
$ws2.Range("A57").Value = "2019.............................................................................     ."
$ws2.Range("B57").Value = 254.412
$ws2.Range("C57").Value = 256.903
$ws2.Range("D57").Value = 255.657
$ws2.Range("E57").Value = 2.3
$ws2.Range("F57").Value = 1.8

$ws2.Range("G57").Formula = '=$D$50/D57'
$ws2.Range("G57").NumberFormat = "0.000"

# --- About sheet: bump the "current data year" and source link --------
$ws1.Range("B4").Value = 2019
$ws1.Range("B6").Value = "https://www.bls.gov/cpi/tables/supplemental-files/historical-cpi-u-201912.pdf"
